$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The image-path cells in column C previously stored paths prefixed with
# "assets/" (e.g. "assets/Course1.png"). Update them to store just the
# bare file name (e.g. "Course1.png").
$ws.Range("C1").Value = "Course1.png"
$ws.Range("C2").Value = "Course2.png"
$ws.Range("C3").Value = "Course3.png"
$ws.Range("C4").Value = "Course4.png"

# Update the view state: scroll the window so row 4 is visible at the top
# and select F4 as the active cell.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F4").Select()
